# DomesticRotterdamDatabase.xlsx - increase hybrid/electric surcharges (transport calculator prices)
# and update the active sheet view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update PRICE column (D) for every POL/vehicle combination. ---
# Each "city block" of 5 rows follows the same vehicle order:
#   CAR, SUV, LARGE SUV, PICKUP, MOTORCYCLE
# and takes the same new surcharge values.
$newPrices = @(920, 1025, 1125, 1230, 715)

# Rows 2-36: destinations Rotterdam (New York, Savannah, Miami, Houston,
# Indianapolis, Los Angeles, San Francisco) - 7 blocks of 5 rows.
$row = 2
for ($block = 0; $block -lt 7; $block++) {
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item($row, 4).Value = $newPrices[$i]
        $row++
    }
}

# --- Update the sheet view: scroll back to the top-left and change the
#     active selection from G45 to D32:D36. ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D32:D36").Select()
